# Adds a new "Método de Asignación" column (K) to the LOCALES report sheet.
# Mirrors the formatting of the existing last header column (J), widens the
# title merge to cover the new column, and sizes the new column like Excel's
# "AutoFit" would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the neighboring header cells onto the new column K
# (a plain Style assignment doesn't stick on a brand-new cell here, so use a
# formats-only paste instead, same end result as picking it up with the
# format painter in the UI).
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K2").Value = "Método de Asignación"

# Extend the title merge across the new column.
$ws.Range("A1:K1").Merge()

# Size column K the way it ends up after an autosize-to-content.
$ws.Columns.Item(11).ColumnWidth = 19.2

# Leave the selection on A2, matching the post-edit workbook state.
$ws.Range("A2").Select()
